$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new transaction was recorded above the most-recent "corporate internet
# share" entry, so insert a fresh row at row 35 (pushing the existing
# history for this sheet down by one) and fill in its September
# details/date cells.
$ws.Rows.Item(35).Insert()

$ws.Range("R35").Value = "corporate internet share"
$ws.Range("S35").Value = "2024-09-09 11:39:30"
